# Apply cryptos list update (Sat Jun  3 08:41:33 UTC 2023 GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column stores values as text (e.g. "27.171.93", "1.0000") rather than
# numbers, so force a Text number format before writing to avoid Excel silently
# reinterpreting values like "1.0000" or "10.17" as numeric.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "27.171.93"
$ws.Range("E2").Value = "  +0.28%  "

# Row 3
$ws.Range("D3").Value = "1.903.36"
$ws.Range("E3").Value = "  +0.83%  "

# Row 4
$ws.Range("E4").Value = "  +0.18%  "

# Row 5
$ws.Range("D5").Value = "306.00"
$ws.Range("E5").Value = "  -0.43%  "

# Row 6
$ws.Range("D6").Value = "1.0000"
$ws.Range("E6").Value = "  +0.10%  "

# Row 7
$ws.Range("D7").Value = "0.5238"
$ws.Range("E7").Value = "  +2.08%  "

# Row 8
$ws.Range("E8").Value = "  +1.11%  "

# Row 9
$ws.Range("D9").Value = "0.07251"

# Row 10
$ws.Range("D10").Value = "21.13"
$ws.Range("E10").Value = "  +0.34%  "

# Row 11
$ws.Range("D11").Value = "0.9006"
$ws.Range("E11").Value = "  -0.44%  "

# Row 12
$ws.Range("D12").Value = "0.08486"
$ws.Range("E12").Value = "  +11.31%  "

# Row 13
$ws.Range("D13").Value = "1.902.66"
$ws.Range("E13").Value = "  +0.75%  "

# Row 14
$ws.Range("D14").Value = "95.07"
$ws.Range("E14").Value = "  +0.66%  "

# Row 15
$ws.Range("D15").Value = "5.288"

# Row 16
$ws.Range("E16").Value = "  +0.17%  "

# Row 17
$ws.Range("D17").Value = "0.000008627"
$ws.Range("E17").Value = "  +1.32%  "

# Row 18
$ws.Range("E18").Value = "  +1.37%  "

# Row 19
$ws.Range("D19").Value = "1.0000"
$ws.Range("E19").Value = "  +0.17%  "

# Row 20
$ws.Range("D20").Value = "27.207.96"
$ws.Range("E20").Value = "  +0.38%  "

# Row 21
$ws.Range("D21").Value = "5.068"
$ws.Range("E21").Value = "  +0.07%  "

# Row 22
$ws.Range("D22").Value = "2.142.80"
$ws.Range("E22").Value = "  +0.36%  "

# Row 23
$ws.Range("D23").Value = "10.60"
$ws.Range("E23").Value = "  +0.53%  "

# Row 24
$ws.Range("D24").Value = "6.424"
$ws.Range("E24").Value = "  +0.23%  "

# Row 25
$ws.Range("D25").Value = "147.27"
$ws.Range("E25").Value = "  +0.62%  "

# Row 26
$ws.Range("E26").Value = "  +4.78%  "

# Row 27
$ws.Range("D27").Value = "1.750"
$ws.Range("E27").Value = "  -2.34%  "

# Row 28
$ws.Range("E28").Value = "  +0.96%  "

# Row 29
$ws.Range("D29").Value = "114.93"
$ws.Range("E29").Value = "  +0.34%  "

# Row 30
$ws.Range("D30").Value = "4.814"
$ws.Range("E30").Value = "  -0.73%  "

# Row 31
$ws.Range("D31").Value = "4.888"
$ws.Range("E31").Value = "  -2.02%  "

# Row 32
$ws.Range("D32").Value = "0.09254"
$ws.Range("E32").Value = "  +0.56%  "

# Row 33
$ws.Range("D33").Value = "0.8079"
$ws.Range("E33").Value = "  +5.03%  "

# Row 34
$ws.Range("D34").Value = "0.05063"
$ws.Range("E34").Value = "  +0.01%  "

# Row 35
$ws.Range("D35").Value = "1.237"
$ws.Range("E35").Value = "  +2.90%  "

# Row 36
$ws.Range("D36").Value = "3.427"
$ws.Range("E36").Value = "  +4.69%  "

# Row 37
$ws.Range("D37").Value = "2.945"
$ws.Range("E37").Value = "  -1.01%  "

# Row 38
$ws.Range("D38").Value = "2.618"
$ws.Range("E38").Value = "  +1.15%  "

# Row 39
$ws.Range("D39").Value = "0.5726"
$ws.Range("E39").Value = "  +1.75%  "

# Row 41
$ws.Range("E41").Value = "  -0.08%  "

# Row 42
$ws.Range("D42").Value = "9.021"
$ws.Range("E42").Value = "  +0.94%  "

# Row 43
$ws.Range("D43").Value = "6.637"
$ws.Range("E43").Value = "  +0.63%  "

# Row 44
$ws.Range("D44").Value = "116.48"
$ws.Range("E44").Value = "  -1.85%  "

# Row 45
$ws.Range("D45").Value = "0.1513"
$ws.Range("E45").Value = "  +0.62%  "

# Row 46
$ws.Range("D46").Value = "0.4862"
$ws.Range("E46").Value = "  +0.97%  "

# Row 47
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "10.17"
$ws.Range("E47").Value = "  +0.64%  "

# Row 48
$ws.Range("B48").Value = "PaxDollar"
$ws.Range("C48").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D48").Value = "0.9998"
$ws.Range("E48").Value = "  +0.08%  "

# Row 49
$ws.Range("D49").Value = "1.615"
$ws.Range("E49").Value = "  +1.34%  "

# Row 50
$ws.Range("D50").Value = "37.46"
$ws.Range("E50").Value = "  +0.65%  "

# Row 51
$ws.Range("D51").Value = "63.90"
$ws.Range("E51").Value = "  -0.38%  "

# Restore the default "General" number format on the Price column now that the
# text values have been written (keeps styling consistent with the rest of the sheet).
$ws.Range("D2:D51").NumberFormat = "General"
